$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 19:47"

# Row 4
$ws.Range("B4").Value = 8679964
$ws.Range("C4").Value = 18313
$ws.Range("D4").Value = 5664735
$ws.Range("E4").Value = 2786603
$ws.Range("G4").Value = 245
$ws.Range("H4").Value = 228626

# Row 5
$ws.Range("B5").Value = 7810405
$ws.Range("C5").Value = 50765
$ws.Range("D5").Value = 7009382
$ws.Range("E5").Value = 683115
$ws.Range("G5").Value = 572
$ws.Range("H5").Value = 117908

# Row 8
$ws.Range("B8").Value = 1110372
$ws.Range("C8").Value = 19851
$ws.Range("G8").Value = 231
$ws.Range("H8").Value = 34752

# Row 33
$ws.Range("B33").Value = 210881
$ws.Range("C33").Value = 1733
$ws.Range("D33").Value = 177307
$ws.Range("E33").Value = 23691

# Row 35
$ws.Range("B35").Value = 190416
$ws.Range("C35").Value = 3685
$ws.Range("D35").Value = 157175
$ws.Range("E35").Value = 30036
$ws.Range("G35").Value = 73
$ws.Range("H35").Value = 3205

# Row 68
$ws.Range("A68").Value = "Irlanda"
$ws.Range("B68").Value = 55261
$ws.Range("C68").Value = 785
$ws.Range("D68").Value = 23364
$ws.Range("E68").Value = 30019
$ws.Range("G68").Value = 7
$ws.Range("H68").Value = 1878

# Row 69
$ws.Range("A69").Value = "Kirguistan"
$ws.Range("B69").Value = 54588
$ws.Range("C69").Value = 582
$ws.Range("D69").Value = 47050
$ws.Range("E69").Value = 6412
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 1126

# Row 73
$ws.Range("A73").Value = "Jordania"
$ws.Range("B73").Value = 48930
$ws.Range("C73").Value = 2489
$ws.Range("D73").Value = 7449
$ws.Range("E73").Value = 40973
$ws.Range("G73").Value = 27
$ws.Range("H73").Value = 508

# Row 74
$ws.Range("A74").Value = "Azerbaiyan"
$ws.Range("B74").Value = 48221
$ws.Range("C74").Value = 803
$ws.Range("D74").Value = 40831
$ws.Range("E74").Value = 6734
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 656

# Row 75
$ws.Range("A75").Value = "Kenia"
$ws.Range("B75").Value = 47843
$ws.Range("C75").Value = 631
$ws.Range("D75").Value = 33421
$ws.Range("E75").Value = 13538
$ws.Range("G75").Value = 14
$ws.Range("H75").Value = 884

# Row 76
$ws.Range("A76").Value = "Ghana"
$ws.Range("B76").Value = 47601
$ws.Range("C76").Value = 63
$ws.Range("D76").Value = 46824
$ws.Range("E76").Value = 463
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 314

# Row 77
$ws.Range("A77").Value = "Tunez"
$ws.Range("B77").Value = 47214
$ws.Range("C77").Value = 1322
$ws.Range("D77").Value = 5032
$ws.Range("E77").Value = 41398
$ws.Range("G77").Value = 44
$ws.Range("H77").Value = 784

# Row 106
$ws.Range("A106").Value = "Mozambique"
$ws.Range("B106").Value = 11748
$ws.Range("C106").Value = 189
$ws.Range("D106").Value = 9234
$ws.Range("E106").Value = 2432
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 82

# Row 107
$ws.Range("A107").Value = "Guinea"
$ws.Range("B107").Value = 11635
$ws.Range("D107").Value = 10474
$ws.Range("E107").Value = 1090
$ws.Range("H107").Value = 71

# Row 122
$ws.Range("B122").Value = 7153
$ws.Range("C122").Value = 866
$ws.Range("E122").Value = 3495

# Row 183
$ws.Range("B183").Value = 461
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 405
$ws.Range("E183").Value = 56

# Row 216
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# Row 217
$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
